$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R4: 2021 header, same formatting as the existing P4/Q4 year headers ---
$q4 = $ws.Range("Q4")
$r4 = $ws.Range("R4")
$r4.Value = 2021
$q4.Copy()
$r4.PasteSpecial(-4122)   # xlPasteFormats - reuse Q4's existing style (no new style created)

# --- R5: new data point, same base formatting as Q5 but with a "0.0" number format ---
$q5 = $ws.Range("Q5")
$r5 = $ws.Range("R5")
$r5.Value = 102.20441221981518
$q5.Copy()
$r5.PasteSpecial(-4122)   # xlPasteFormats - start from Q5's font/border/alignment
$r5.NumberFormat = "0.0"  # then layer on the new "0.0" number format (numFmtId 167)

$excel.CutCopyMode = 0

# --- restore the active selection to match the saved workbook state ---
$null = $ws.Range("S9").Select()
